$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet (Loan RBI, Variable Instalments): insert a new
# blank column before the existing "Late" column so the schedule gains a
# spacer column between "In Advance" and "Late" (old N/O/P -> new O/P/Q).
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab (was "Input"), with K19
# selected.
$wsSchedule.Activate()
$wsSchedule.Range("K19").Select()
